# Update cryptocurrency price/volume figures in the tracking sheet.
# Mirrors the "Updated symbol list" GitHub Actions commit: refreshes the
# Price (column D) and Volume(1h) (column E) text values for each coin row.

function Set-TextValue($ws, $row, $col, $val) {
    # Assign as literal text so numeric-looking strings (prices, percentages)
    # keep their exact formatting (trailing zeros, "%", thousands separators)
    # instead of being auto-coerced into numbers by the Value setter.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 2 4 '246.26'
Set-TextValue $ws 3 4 '29.54'
Set-TextValue $ws 3 5 '8.18%'
Set-TextValue $ws 4 4 '5.190'
Set-TextValue $ws 4 5 '2.93%'
Set-TextValue $ws 5 4 '0.05701'
Set-TextValue $ws 5 5 '0.44%'
Set-TextValue $ws 6 4 '6.591'
Set-TextValue $ws 6 5 '1.84%'
Set-TextValue $ws 7 4 '0.8580'
Set-TextValue $ws 7 5 '4.28%'
Set-TextValue $ws 8 4 '0.8777'
Set-TextValue $ws 8 5 '4.04%'
Set-TextValue $ws 9 4 '0.1371'
Set-TextValue $ws 9 5 '3.20%'
Set-TextValue $ws 10 4 '0.07109'
Set-TextValue $ws 10 5 '2.68%'
Set-TextValue $ws 11 4 '0.02871'
Set-TextValue $ws 11 5 '-0.45%'
Set-TextValue $ws 12 4 '0.09387'
Set-TextValue $ws 12 5 '-0.04%'
Set-TextValue $ws 13 4 '0.001526'
Set-TextValue $ws 13 5 '1.15%'
Set-TextValue $ws 14 4 '0.04156'
Set-TextValue $ws 14 5 '0.46%'
Set-TextValue $ws 15 4 '0.0005981'
Set-TextValue $ws 15 5 '-0.48%'
Set-TextValue $ws 16 4 '0.006127'
Set-TextValue $ws 16 5 '-0.87%'
Set-TextValue $ws 17 5 '5,107.96%'
Set-TextValue $ws 18 5 '-0.88%'
Set-TextValue $ws 19 4 '3.057'
Set-TextValue $ws 19 5 '1.91%'
Set-TextValue $ws 20 4 '2.179'
Set-TextValue $ws 20 5 '-5.61%'
Set-TextValue $ws 21 4 '0.3145'
Set-TextValue $ws 21 5 '1.01%'
Set-TextValue $ws 22 4 '0.03269'
Set-TextValue $ws 22 5 '3.23%'
Set-TextValue $ws 23 5 '3.67%'
Set-TextValue $ws 24 4 '3.472'
Set-TextValue $ws 24 5 '-2.41%'
Set-TextValue $ws 26 5 '31.54%'
Set-TextValue $ws 27 4 '0.001218'
Set-TextValue $ws 27 5 '-0.06%'
Set-TextValue $ws 28 4 '0.0001210'
Set-TextValue $ws 28 5 '23.45%'
Set-TextValue $ws 40 4 '0.03746'
Set-TextValue $ws 40 5 '2.23%'
Set-TextValue $ws 41 4 '0.005679'
Set-TextValue $ws 41 5 '-6.13%'
Set-TextValue $ws 42 4 '0.1073'
Set-TextValue $ws 42 5 '1.88%'
Set-TextValue $ws 43 4 '0.002540'
Set-TextValue $ws 43 5 '11.50%'
Set-TextValue $ws 44 4 '0.009415'
Set-TextValue $ws 44 5 '-11.39%'
Set-TextValue $ws 45 4 '0.00005117'
Set-TextValue $ws 45 5 '-3.64%'
Set-TextValue $ws 46 5 '0.00%'
Set-TextValue $ws 47 4 '0.07101'
Set-TextValue $ws 47 5 '-30.03%'
Set-TextValue $ws 48 4 '0.002667'
Set-TextValue $ws 48 5 '4.13%'
Set-TextValue $ws 49 4 '0.00002100'
Set-TextValue $ws 49 5 '0.00%'
Set-TextValue $ws 50 4 '0.0002000'
Set-TextValue $ws 50 5 '0.00%'
